$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("YDS")

$ws.Range("B2").Value = "NCT(2.8897424023070566, 1.1396013648214876, 0.24127482754312501, 2.606802356152179)"
$ws.Range("C2").Value = "NIG(1.4179851095053033, 1.0181052132002812, 3.8869359883094967, 6.283390878215604)"
$ws.Range("D2").Value = "NIG(0.6908286094213698, 0.4555737041387318, 1.5878312762088909, 2.800272387263946)"
$ws.Range("E2").Value = "JSU(-1.4944948884536071, 1.3426110106194833, 1.738018519279759, 5.302845331089905)"
